# This script inserts two new data rows at row 223 of Sheet1 (pushing the
# existing rows 223-275 down to 225-277), and populates the two new rows
# with the values required by the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 223/224; everything below shifts down.
$ws.Rows("223:224").Insert()

# ---- New row 223 ----
$ws.Cells.Item(223, 1).Value  = 5
$ws.Cells.Item(223, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(223, 3).Value  = "Maule"
$ws.Cells.Item(223, 4).Value  = "1/6/2023"
$ws.Cells.Item(223, 5).Value  = 7
$ws.Cells.Item(223, 6).Value  = 100112021
$ws.Cells.Item(223, 7).Value  = "Ají"
$ws.Cells.Item(223, 8).Value  = "Americana (o)"
$ws.Cells.Item(223, 9).Value  = "Primera"
$ws.Cells.Item(223, 10).Value = 150
$ws.Cells.Item(223, 11).Value = 12000
$ws.Cells.Item(223, 12).Value = 12000
$ws.Cells.Item(223, 13).Value = 12000
$ws.Cells.Item(223, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(223, 15).Value = "Región del Maule"
$ws.Cells.Item(223, 16).Value = 800
$ws.Cells.Item(223, 17).Value = 15
$ws.Cells.Item(223, 18).Value = "Hortaliza"

# ---- New row 224 ----
$ws.Cells.Item(224, 1).Value  = 5
$ws.Cells.Item(224, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(224, 3).Value  = "Maule"
$ws.Cells.Item(224, 4).Value  = "1/6/2023"
$ws.Cells.Item(224, 5).Value  = 7
$ws.Cells.Item(224, 6).Value  = 100112021
$ws.Cells.Item(224, 7).Value  = "Ají"
$ws.Cells.Item(224, 8).Value  = "Cacho cabra verde"
$ws.Cells.Item(224, 9).Value  = "Primera"
$ws.Cells.Item(224, 10).Value = 80
$ws.Cells.Item(224, 11).Value = 12000
$ws.Cells.Item(224, 12).Value = 12000
$ws.Cells.Item(224, 13).Value = 12000
$ws.Cells.Item(224, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(224, 15).Value = "Región del Maule"
$ws.Cells.Item(224, 16).Value = 800
$ws.Cells.Item(224, 17).Value = 15
$ws.Cells.Item(224, 18).Value = "Hortaliza"
